$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1067.3334
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -85
$ws.Range("H41").Value = 366.33334
$ws.Range("I41").Value = 349.5
$ws.Range("K41").Value = 349.5
$ws.Range("M41").Value = 90.5
$ws.Range("H43").Value = 1446.6666
$ws.Range("I43").Value = 1588
$ws.Range("J43").Value = 740
$ws.Range("K43").Value = 1588
$ws.Range("L43").Value = 740
$ws.Range("M43").Value = -1519
$ws.Range("N43").Value = -878
$ws.Range("H53").Value = 588.6667
$ws.Range("I53").Value = 486
$ws.Range("J53").Value = 794
$ws.Range("K53").Value = 486
$ws.Range("L53").Value = 794
$ws.Range("M53").Value = 151
$ws.Range("N53").Value = -2068
$ws.Range("H112").Value = 1499.25
$ws.Range("J112").Value = 1499.25
$ws.Range("L112").Value = 4497.75
$ws.Range("N112").Value = -6713.75
$ws.Range("H121").Value = 1909
$ws.Range("J121").Value = 1909
$ws.Range("L121").Value = 5727
$ws.Range("N121").Value = -9221
$ws.Range("H125").Value = 8492.666999999999
$ws.Range("I125").Value = 8243
$ws.Range("K125").Value = 74187
$ws.Range("M125").Value = -71727
$ws.Range("H131").Value = 9999
$ws.Range("I131").Value = 9999
$ws.Range("K131").Value = 29997
$ws.Range("M131").Value = -24957
$ws.Range("H135").Value = 2275.625
$ws.Range("I135").Value = 1460
$ws.Range("J135").Value = 2765
$ws.Range("K135").Value = 13140
$ws.Range("L135").Value = 24885
$ws.Range("M135").Value = -10605
$ws.Range("N135").Value = -29955
$ws.Range("H138").Value = 5216.0786
$ws.Range("I138").Value = 1467.8
$ws.Range("J138").Value = 6130.2925
$ws.Range("K138").Value = 4403.4
$ws.Range("L138").Value = 18390.8775
$ws.Range("M138").Value = 736.6000000000004
$ws.Range("N138").Value = -28670.8775
$ws.Range("H141").Value = 6747.5
$ws.Range("I141").Value = 8330
$ws.Range("K141").Value = 24990
$ws.Range("M141").Value = -19810
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5094.5
$ws.Range("J61").Value = 3998.5
$ws.Range("L61").Value = 3998.5
$ws.Range("N61").Value = -4422.5
$ws.Range("H132").Value = 5177.6
$ws.Range("I132").Value = 4890
$ws.Range("J132").Value = 5249.5
$ws.Range("K132").Value = 14670
$ws.Range("L132").Value = 15748.5
$ws.Range("M132").Value = -12140
$ws.Range("N132").Value = -20808.5
$ws.Range("H135").Value = 133333
$ws.Range("J135").Value = 133333
$ws.Range("L135").Value = 133333
$ws.Range("N135").Value = -143473
$ws.Range("H136").Value = 5094.5
$ws.Range("J136").Value = 3998.5
$ws.Range("L136").Value = 11995.5
$ws.Range("N136").Value = -17095.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3243.0625
$ws.Range("I20").Value = 3422.6667
$ws.Range("K20").Value = 3422.6667
$ws.Range("M20").Value = -3175.6667
$ws.Range("H94").Value = 1060.75
$ws.Range("I94").Value = 821.75
$ws.Range("K94").Value = 821.75
$ws.Range("M94").Value = -370.75
$ws.Range("H134").Value = 2045.3158
$ws.Range("I134").Value = 1750.75
$ws.Range("J134").Value = 3616.3333
$ws.Range("K134").Value = 5252.25
$ws.Range("L134").Value = 10848.9999
$ws.Range("M134").Value = -2717.25
$ws.Range("N134").Value = -15918.9999
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5950.143
$ws.Range("I31").Value = 2926.5
$ws.Range("K31").Value = 2926.5
$ws.Range("M31").Value = -2631.5
$ws.Range("H34").Value = 5950.143
$ws.Range("I34").Value = 2926.5
$ws.Range("K34").Value = 2926.5
$ws.Range("M34").Value = -2724.5
$ws.Range("H99").Value = 3697.5
$ws.Range("I99").Value = 3666
$ws.Range("J99").Value = 3729
$ws.Range("K99").Value = 3666
$ws.Range("L99").Value = 3729
$ws.Range("M99").Value = -2168
$ws.Range("N99").Value = -6725
$ws.Range("H126").Value = 3697.5
$ws.Range("I126").Value = 3666
$ws.Range("J126").Value = 3729
$ws.Range("K126").Value = 10998
$ws.Range("L126").Value = 11187
$ws.Range("M126").Value = -8528
$ws.Range("N126").Value = -16127
$ws.Range("H134").Value = 4013
$ws.Range("I134").Value = 4013
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 12039
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -9504
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6012.625
$ws.Range("I80").Value = 4419.2
$ws.Range("J80").Value = 8668.333000000001
$ws.Range("K80").Value = 4419.2
$ws.Range("L80").Value = 8668.333000000001
$ws.Range("M80").Value = -3421.2
$ws.Range("N80").Value = -10664.333
$ws.Range("H83").Value = 6012.625
$ws.Range("I83").Value = 4419.2
$ws.Range("J83").Value = 8668.333000000001
$ws.Range("K83").Value = 22096
$ws.Range("L83").Value = 43341.665
$ws.Range("M83").Value = -17104
$ws.Range("N83").Value = -53325.665
$ws.Range("H102").Value = 3580.4119
$ws.Range("I102").Value = 3580.4119
$ws.Range("K102").Value = 3580.4119
$ws.Range("M102").Value = -1958.4119
$ws.Range("H132").Value = 6663.3335
$ws.Range("J132").Value = 6663.3335
$ws.Range("L132").Value = 19990.0005
$ws.Range("N132").Value = -25050.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H82").Value = 1495
$ws.Range("I82").Value = 1490
$ws.Range("K82").Value = 1490
$ws.Range("M82").Value = -1129
$ws.Range("H85").Value = 1495
$ws.Range("I85").Value = 1490
$ws.Range("K85").Value = 1490
$ws.Range("M85").Value = -242
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2821.6123
$ws.Range("I132").Value = 2174.9714
$ws.Range("K132").Value = 6524.914199999999
$ws.Range("M132").Value = -3994.914199999999
$ws.Range("H136").Value = 9928.933999999999
$ws.Range("I136").Value = 10665.538
$ws.Range("K136").Value = 31996.614
$ws.Range("M136").Value = -29446.614
